# Scheduled runner update: refresh cached market-price / profit figures
# (currentAveragePrice*, LevePriceNQ/HQ, LeveProfitNQ/HQ columns H-N)
# across the per-job Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1831.1765
$ws.Range("I28").Value = 1831.1765
$ws.Range("K28").Value = 1831.1765
$ws.Range("M28").Value = -1346.1765
$ws.Range("H74").Value = 3312.875
$ws.Range("J74").Value = 3500
$ws.Range("L74").Value = 3500
$ws.Range("N74").Value = -5372
$ws.Range("H77").Value = 3312.875
$ws.Range("J77").Value = 3500
$ws.Range("L77").Value = 17500
$ws.Range("N77").Value = -26860
$ws.Range("H92").Value = 1539.1666
$ws.Range("I92").Value = 1291.3636
$ws.Range("J92").Value = 1928.5714
$ws.Range("K92").Value = 1291.3636
$ws.Range("L92").Value = 1928.5714
$ws.Range("M92").Value = -43.36359999999991
$ws.Range("N92").Value = -4424.5714
$ws.Range("H94").Value = 1993.3334
$ws.Range("I94").Value = 1993.3334
$ws.Range("K94").Value = 1993.3334
$ws.Range("M94").Value = -1542.3334
$ws.Range("H100").Value = 1490.7693
$ws.Range("I100").Value = 1366.4445
$ws.Range("K100").Value = 1366.4445
$ws.Range("M100").Value = -825.4445000000001
$ws.Range("H103").Value = 1174.8334
$ws.Range("I103").Value = 589.6
$ws.Range("J103").Value = 1592.8572
$ws.Range("K103").Value = 1768.8
$ws.Range("L103").Value = 4778.571599999999
$ws.Range("M103").Value = -1182.8
$ws.Range("N103").Value = -5950.571599999999
$ws.Range("H106").Value = 9451.286
$ws.Range("I106").Value = 9947.538
$ws.Range("K106").Value = 9947.538
$ws.Range("M106").Value = -9316.538
$ws.Range("H107").Value = 2765.8235
$ws.Range("I107").Value = 2332.1
$ws.Range("K107").Value = 2332.1
$ws.Range("M107").Value = -412.0999999999999
$ws.Range("H137").Value = 1366.0975
$ws.Range("I137").Value = 1016.4231
$ws.Range("K137").Value = 3049.2693
$ws.Range("M137").Value = -499.2692999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 20836044
$ws.Range("I102").Value = 27780212
$ws.Range("K102").Value = 27780212
$ws.Range("M102").Value = -27778590
$ws.Range("H110").Value = 998.8889
$ws.Range("I110").Value = 606.8261
$ws.Range("J110").Value = 3253.25
$ws.Range("K110").Value = 606.8261
$ws.Range("L110").Value = 3253.25
$ws.Range("M110").Value = 1438.1739
$ws.Range("N110").Value = -7343.25
$ws.Range("H122").Value = 1415.9166
$ws.Range("I122").Value = 1071.1333
$ws.Range("J122").Value = 1990.5555
$ws.Range("K122").Value = 3213.3999
$ws.Range("L122").Value = 5971.666499999999
$ws.Range("M122").Value = -763.3998999999999
$ws.Range("N122").Value = -10871.6665
$ws.Range("H132").Value = 3836.6667
$ws.Range("I132").Value = 4302.778
$ws.Range("J132").Value = 3370.5557
$ws.Range("K132").Value = 12908.334
$ws.Range("L132").Value = 10111.6671
$ws.Range("M132").Value = -10378.334
$ws.Range("N132").Value = -15171.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 37038030
$ws.Range("I99").Value = 45455548
$ws.Range("J99").Value = 976.2
$ws.Range("K99").Value = 45455548
$ws.Range("L99").Value = 976.2
$ws.Range("M99").Value = -45454050
$ws.Range("N99").Value = -3972.2
$ws.Range("H105").Value = 45455436
$ws.Range("I105").Value = 45455436
$ws.Range("K105").Value = 45455436
$ws.Range("M105").Value = -45453689
$ws.Range("H107").Value = 1393.3846
$ws.Range("I107").Value = 950.125
$ws.Range("K107").Value = 950.125
$ws.Range("M107").Value = 969.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 314.54544
$ws.Range("I7").Value = 133.6
$ws.Range("K7").Value = 133.6
$ws.Range("M7").Value = -20.59999999999999
$ws.Range("H22").Value = 100227.43
$ws.Range("I22").Value = 197.75
$ws.Range("K22").Value = 197.75
$ws.Range("M22").Value = 152.25
$ws.Range("H31").Value = 1487.2667
$ws.Range("I31").Value = 1758
$ws.Range("J31").Value = 1351.9
$ws.Range("K31").Value = 1758
$ws.Range("L31").Value = 1351.9
$ws.Range("M31").Value = -1463
$ws.Range("N31").Value = -1941.9
$ws.Range("H34").Value = 1487.2667
$ws.Range("I34").Value = 1758
$ws.Range("J34").Value = 1351.9
$ws.Range("K34").Value = 1758
$ws.Range("L34").Value = 1351.9
$ws.Range("M34").Value = -1556
$ws.Range("N34").Value = -1755.9
$ws.Range("H86").Value = 3196751.5
$ws.Range("I86").Value = 8348707.5
$ws.Range("J86").Value = 26317
$ws.Range("K86").Value = 8348707.5
$ws.Range("L86").Value = 26317
$ws.Range("M86").Value = -8347584.5
$ws.Range("N86").Value = -28563
$ws.Range("H89").Value = 3196751.5
$ws.Range("I89").Value = 8348707.5
$ws.Range("J89").Value = 26317
$ws.Range("K89").Value = 41743537.5
$ws.Range("L89").Value = 131585
$ws.Range("M89").Value = -41737921.5
$ws.Range("N89").Value = -142817
$ws.Range("H105").Value = 694.2857
$ws.Range("I105").Value = 676.6667
$ws.Range("K105").Value = 676.6667
$ws.Range("M105").Value = 1070.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 500.07144
$ws.Range("I7").Value = 545.5454999999999
$ws.Range("K7").Value = 1636.6365
$ws.Range("M7").Value = -1524.6365
$ws.Range("H40").Value = 237
$ws.Range("I40").Value = 102.111115
$ws.Range("J40").Value = 317.93332
$ws.Range("K40").Value = 408.44446
$ws.Range("L40").Value = 1271.73328
$ws.Range("M40").Value = -339.44446
$ws.Range("N40").Value = -1409.73328

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 3093
$ws.Range("J54").Value = 3093
$ws.Range("L54").Value = 3093
$ws.Range("N54").Value = -3873
$ws.Range("H97").Value = 749.5833
$ws.Range("I97").Value = 799.2857
$ws.Range("K97").Value = 799.2857
$ws.Range("M97").Value = -303.2857
$ws.Range("H107").Value = 672.4545000000001
$ws.Range("I107").Value = 826
$ws.Range("J107").Value = 584.7143
$ws.Range("K107").Value = 826
$ws.Range("L107").Value = 584.7143
$ws.Range("M107").Value = 1094
$ws.Range("N107").Value = -4424.7143
$ws.Range("H113").Value = 1345.4
$ws.Range("I113").Value = 1284.25
$ws.Range("J113").Value = 1590
$ws.Range("K113").Value = 1284.25
$ws.Range("L113").Value = 1590
$ws.Range("M113").Value = 885.75
$ws.Range("N113").Value = -5930
$ws.Range("H122").Value = 2186.9524
$ws.Range("I122").Value = 2200.5715
$ws.Range("J122").Value = 2159.7144
$ws.Range("K122").Value = 6601.7145
$ws.Range("L122").Value = 6479.1432
$ws.Range("M122").Value = -4151.7145
$ws.Range("N122").Value = -11379.1432
$ws.Range("H132").Value = 3131.348
$ws.Range("I132").Value = 3056.7856
$ws.Range("K132").Value = 9170.356800000001
$ws.Range("M132").Value = -6640.356800000001
$ws.Range("H134").Value = 27491.2
$ws.Range("J134").Value = 27491.2
$ws.Range("L134").Value = 82473.60000000001
$ws.Range("N134").Value = -87543.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1048.25
$ws.Range("I93").Value = 897.6667
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 897.6667
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 350.3333
$ws.Range("N93").Value = -3996
$ws.Range("H100").Value = 1623.8572
$ws.Range("I100").Value = 1392.3334
$ws.Range("J100").Value = 1797.5
$ws.Range("K100").Value = 1392.3334
$ws.Range("L100").Value = 1797.5
$ws.Range("M100").Value = -851.3334
$ws.Range("N100").Value = -2879.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 212.5
$ws.Range("I81").Value = 212.5
$ws.Range("K81").Value = 425
$ws.Range("M81").Value = 636
$ws.Range("H84").Value = 212.5
$ws.Range("I84").Value = 212.5
$ws.Range("K84").Value = 2125
$ws.Range("M84").Value = 3179
$ws.Range("H96").Value = 2489.3845
$ws.Range("I96").Value = 1894.75
$ws.Range("J96").Value = 3440.8
$ws.Range("K96").Value = 1894.75
$ws.Range("L96").Value = 3440.8
$ws.Range("M96").Value = -521.75
$ws.Range("N96").Value = -6186.8
